# Applies the "Fixing network data cleanining scripts" edit to SOUTH_DAKOTA_2016.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) to the cleaned/normalized column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Normalize capitalization of connector words ("de" -> "De", "del" -> "Del",
#    "el" -> "El", "la" -> "La", "los" -> "Los") inside specific place names.
$ws.Range("B5").Value   = "Rincón De Romos"
$ws.Range("A24").Value  = "Ciudad De México"
$ws.Range("A32").Value  = "Estado De México"
$ws.Range("B32").Value  = "Ecatepec De Morelos"
$ws.Range("B33").Value  = "Naucalpan De Juárez"
$ws.Range("B44").Value  = "Purísima Del Rincón"
$ws.Range("B47").Value  = "Santa Cruz De Juventino Rosas"
$ws.Range("B50").Value  = "Atenango Del Río"
$ws.Range("B52").Value  = "Coyuca De Catalán"
$ws.Range("B55").Value  = "Huitzuco De Los Figueroa"
$ws.Range("B58").Value  = "Técpan De Galeana"
$ws.Range("B63").Value  = "Pachuca De Soto"
$ws.Range("B64").Value  = "Tulancingo De Bravo"
$ws.Range("B67").Value  = "Atotonilco El Alto"
$ws.Range("B74").Value  = "San Diego De Alejandría"
$ws.Range("B78").Value  = "Tepatitlán De Morelos"
$ws.Range("B79").Value  = "Tlajomulco De Zúñiga"
$ws.Range("B84").Value  = "Zapotlán El Grande"
$ws.Range("B112").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B113").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B114").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B117").Value = "San Juan Del Río"
$ws.Range("B119").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B121").Value = "Zapotitlán Del Río"
$ws.Range("B126").Value = "Cuetzalan Del Progreso"
$ws.Range("B127").Value = "Huitzilan De Serdán"
$ws.Range("B128").Value = "Izúcar De Matamoros"
$ws.Range("B134").Value = "Landa De Matamoros"
$ws.Range("B139").Value = "Mexquitic De Carmona"
$ws.Range("B142").Value = "Tanquián De Escobedo"
$ws.Range("B161").Value = "Cosamaloapan De Carpio"
$ws.Range("B162").Value = "Ignacio De La Llave"
$ws.Range("B166").Value = "Martínez De La Torre"
$ws.Range("B169").Value = "Soledad De Doblado"

# 3. Remove the trailing footnote / metadata rows (formerly rows 178-182 and
#    476-480) so the sheet's real data now ends at row 176. Deleting the
#    entire row range also shrinks the worksheet dimension to A1:D176.
$ws.Range("A177:A1048576").EntireRow.Delete()
